$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 70000
$ws.Range("B4").Value = 1071428.571428571

$ws.Range("B14").Value = 714285.7142857143

$ws.Range("B24").Value = 1071428.571428571

$ws.Range("B32").Value = 1141428.571428571
$ws.Range("B33").Value = 714285.7142857143
$ws.Range("B34").Value = 1071428.571428571

$ws.Range("A35").Value = "Tổng lương tại HỆ THỐNG"
$ws.Range("B35").Value = 2927142.857142857
